# Auto-generated edit script: updates computed price/profit columns (H:N)
# on each Leve sheet to reflect refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 189.76471
$ws.Range("I9").Value = 146.23077
$ws.Range("J9").Value = 331.25
$ws.Range("K9").Value = 146.23077
$ws.Range("L9").Value = 331.25
$ws.Range("M9").Value = 22.76922999999999
$ws.Range("N9").Value = -669.25
$ws.Range("H12").Value = 1075
$ws.Range("I12").Value = 1075
$ws.Range("K12").Value = 1075
$ws.Range("M12").Value = -905
$ws.Range("H17").Value = 405.6875
$ws.Range("J17").Value = 405.6875
$ws.Range("L17").Value = 1217.0625
$ws.Range("N17").Value = -1553.0625
$ws.Range("H18").Value = 3349.8462
$ws.Range("I18").Value = 2795.6667
$ws.Range("K18").Value = 2795.6667
$ws.Range("M18").Value = -2511.6667
$ws.Range("H29").Value = 6249.5
$ws.Range("I29").Value = 4500
$ws.Range("K29").Value = 13500
$ws.Range("M29").Value = -13219
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H41").Value = 843.1539
$ws.Range("I41").Value = 388.83334
$ws.Range("J41").Value = 1232.5714
$ws.Range("K41").Value = 388.83334
$ws.Range("L41").Value = 1232.5714
$ws.Range("M41").Value = 51.16665999999998
$ws.Range("N41").Value = -2112.5714
$ws.Range("H43").Value = 3298826
$ws.Range("I43").Value = 4793088
$ws.Range("J43").Value = 11449.8
$ws.Range("K43").Value = 4793088
$ws.Range("L43").Value = 11449.8
$ws.Range("M43").Value = -4793019
$ws.Range("N43").Value = -11587.8
$ws.Range("H46").Value = 294705.72
$ws.Range("J46").Value = 343323.34
$ws.Range("L46").Value = 1029970.02
$ws.Range("N46").Value = -1030208.02
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("J51").Value = 2839.4285
$ws.Range("L51").Value = 2839.4285
$ws.Range("N51").Value = -3807.4285
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H60").Value = 294705.72
$ws.Range("J60").Value = 343323.34
$ws.Range("L60").Value = 1029970.02
$ws.Range("N60").Value = -1030938.02
$ws.Range("H69").Value = 12458.538
$ws.Range("J69").Value = 12692.625
$ws.Range("L69").Value = 38077.875
$ws.Range("N69").Value = -39825.875
$ws.Range("H70").Value = 7594.8667
$ws.Range("I70").Value = 1297
$ws.Range("J70").Value = 8044.7144
$ws.Range("K70").Value = 3891
$ws.Range("L70").Value = 24134.1432
$ws.Range("M70").Value = -3621
$ws.Range("N70").Value = -24674.1432
$ws.Range("H72").Value = 12458.538
$ws.Range("J72").Value = 12692.625
$ws.Range("L72").Value = 114233.625
$ws.Range("N72").Value = -122969.625
$ws.Range("H73").Value = 7594.8667
$ws.Range("I73").Value = 1297
$ws.Range("J73").Value = 8044.7144
$ws.Range("K73").Value = 3891
$ws.Range("L73").Value = 24134.1432
$ws.Range("M73").Value = -2955
$ws.Range("N73").Value = -26006.1432
$ws.Range("H76").Value = 2623.8572
$ws.Range("I76").Value = 2623.8572
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2623.8572
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2308.8572
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 2623.8572
$ws.Range("I79").Value = 2623.8572
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2623.8572
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -1531.8572
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 2780.3635
$ws.Range("I86").Value = 2557.5
$ws.Range("K86").Value = 2557.5
$ws.Range("M86").Value = -1434.5
$ws.Range("H89").Value = 2780.3635
$ws.Range("I89").Value = 2557.5
$ws.Range("K89").Value = 12787.5
$ws.Range("M89").Value = -7171.5
$ws.Range("H107").Value = 933.75
$ws.Range("J107").Value = 1842.8889
$ws.Range("L107").Value = 1842.8889
$ws.Range("N107").Value = -5682.8889
$ws.Range("H116").Value = 4181.091
$ws.Range("I116").Value = 3999
$ws.Range("J116").Value = 4399.6
$ws.Range("K116").Value = 3999
$ws.Range("L116").Value = 4399.6
$ws.Range("M116").Value = -557
$ws.Range("N116").Value = -11283.6
$ws.Range("H125").Value = 2047.6666
$ws.Range("I125").Value = 1968.8
$ws.Range("K125").Value = 17719.2
$ws.Range("M125").Value = -15259.2
$ws.Range("H132").Value = 4160.84
$ws.Range("I132").Value = 4160.84
$ws.Range("K132").Value = 12482.52
$ws.Range("M132").Value = -9952.52
$ws.Range("H135").Value = 2566.75
$ws.Range("I135").Value = 2743.2727
$ws.Range("K135").Value = 24689.4543
$ws.Range("M135").Value = -22154.4543
$ws.Range("H137").Value = 6352.5
$ws.Range("I137").Value = 5069.5713
$ws.Range("J137").Value = 15333
$ws.Range("K137").Value = 15208.7139
$ws.Range("L137").Value = 45999
$ws.Range("M137").Value = -12658.7139
$ws.Range("N137").Value = -51099

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 285.2857
$ws.Range("I5").Value = 316.16666
$ws.Range("K5").Value = 316.16666
$ws.Range("M5").Value = -204.16666
$ws.Range("H32").Value = 2124.1296
$ws.Range("I32").Value = 2225.3262
$ws.Range("J32").Value = 1542.25
$ws.Range("K32").Value = 2225.3262
$ws.Range("L32").Value = 1542.25
$ws.Range("M32").Value = -1938.3262
$ws.Range("N32").Value = -2116.25
$ws.Range("H45").Value = 3585.158
$ws.Range("I45").Value = 3284.3333
$ws.Range("K45").Value = 3284.3333
$ws.Range("M45").Value = -2907.3333
$ws.Range("H74").Value = 9020.799999999999
$ws.Range("I74").Value = 10199.25
$ws.Range("J74").Value = 8235.166999999999
$ws.Range("K74").Value = 10199.25
$ws.Range("L74").Value = 8235.166999999999
$ws.Range("M74").Value = -9325.25
$ws.Range("N74").Value = -9983.166999999999
$ws.Range("H77").Value = 9020.799999999999
$ws.Range("I77").Value = 10199.25
$ws.Range("J77").Value = 8235.166999999999
$ws.Range("K77").Value = 50996.25
$ws.Range("L77").Value = 41175.835
$ws.Range("M77").Value = -46628.25
$ws.Range("N77").Value = -49911.835
$ws.Range("H110").Value = 2189.3
$ws.Range("I110").Value = 2245.3333
$ws.Range("J110").Value = 2105.25
$ws.Range("K110").Value = 2245.3333
$ws.Range("L110").Value = 2105.25
$ws.Range("M110").Value = -200.3332999999998
$ws.Range("N110").Value = -6195.25
$ws.Range("H111").Value = 83999
$ws.Range("J111").Value = 83999
$ws.Range("L111").Value = 83999
$ws.Range("N111").Value = -92179
$ws.Range("H122").Value = 2809.3333
$ws.Range("I122").Value = 2481.9092
$ws.Range("K122").Value = 7445.7276
$ws.Range("M122").Value = -4995.7276
$ws.Range("H132").Value = 6247.8423
$ws.Range("I132").Value = 3713.4
$ws.Range("J132").Value = 9063.888999999999
$ws.Range("K132").Value = 11140.2
$ws.Range("L132").Value = 27191.667
$ws.Range("M132").Value = -8610.200000000001
$ws.Range("N132").Value = -32251.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 285.2857
$ws.Range("I4").Value = 316.16666
$ws.Range("K4").Value = 316.16666
$ws.Range("M4").Value = -201.16666
$ws.Range("H22").Value = 271.8
$ws.Range("I22").Value = 268.66666
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 268.66666
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -95.66665999999998
$ws.Range("N22").Value = -646
$ws.Range("H53").Value = 96899
$ws.Range("J53").Value = 96899
$ws.Range("L53").Value = 96899
$ws.Range("N53").Value = -98047
$ws.Range("H96").Value = 20419.285
$ws.Range("I96").Value = 14563.333
$ws.Range("J96").Value = 55555
$ws.Range("K96").Value = 14563.333
$ws.Range("L96").Value = 55555
$ws.Range("M96").Value = -11817.333
$ws.Range("N96").Value = -61047
$ws.Range("H105").Value = 4423.55
$ws.Range("I105").Value = 3504.7273
$ws.Range("J105").Value = 5546.5557
$ws.Range("K105").Value = 3504.7273
$ws.Range("L105").Value = 5546.5557
$ws.Range("M105").Value = -1757.7273
$ws.Range("N105").Value = -9040.555700000001
$ws.Range("H107").Value = 12949.5
$ws.Range("I107").Value = 14303.333
$ws.Range("J107").Value = 8888
$ws.Range("K107").Value = 14303.333
$ws.Range("L107").Value = 8888
$ws.Range("M107").Value = -12383.333
$ws.Range("N107").Value = -12728

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 567.5
$ws.Range("I22").Value = 551.6875
$ws.Range("J22").Value = 694
$ws.Range("K22").Value = 551.6875
$ws.Range("L22").Value = 694
$ws.Range("M22").Value = -201.6875
$ws.Range("N22").Value = -1394
$ws.Range("H31").Value = 3543.111
$ws.Range("I31").Value = 2915.4546
$ws.Range("J31").Value = 3746.1765
$ws.Range("K31").Value = 2915.4546
$ws.Range("L31").Value = 3746.1765
$ws.Range("M31").Value = -2620.4546
$ws.Range("N31").Value = -4336.1765
$ws.Range("H34").Value = 3543.111
$ws.Range("I34").Value = 2915.4546
$ws.Range("J34").Value = 3746.1765
$ws.Range("K34").Value = 2915.4546
$ws.Range("L34").Value = 3746.1765
$ws.Range("M34").Value = -2713.4546
$ws.Range("N34").Value = -4150.1765
$ws.Range("H52").Value = 149999
$ws.Range("J52").Value = 149999
$ws.Range("L52").Value = 149999
$ws.Range("N52").Value = -150587
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H58").Value = 5936.091
$ws.Range("I58").Value = 4794.4287
$ws.Range("J58").Value = 6468.8667
$ws.Range("K58").Value = 4794.4287
$ws.Range("L58").Value = 6468.8667
$ws.Range("M58").Value = -4591.4287
$ws.Range("N58").Value = -6874.8667
$ws.Range("H99").Value = 5962.56
$ws.Range("I99").Value = 5084.2
$ws.Range("K99").Value = 5084.2
$ws.Range("M99").Value = -3586.2
$ws.Range("H107").Value = 412.28125
$ws.Range("I107").Value = 328.9
$ws.Range("K107").Value = 328.9
$ws.Range("M107").Value = 1591.1
$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 35000
$ws.Range("L111").Value = 35000
$ws.Range("N111").Value = -43180
$ws.Range("H126").Value = 5962.56
$ws.Range("I126").Value = 5084.2
$ws.Range("K126").Value = 15252.6
$ws.Range("M126").Value = -12782.6
$ws.Range("H132").Value = 4139.0586
$ws.Range("I132").Value = 1053
$ws.Range("J132").Value = 5088.615
$ws.Range("K132").Value = 3159
$ws.Range("L132").Value = 15265.845
$ws.Range("M132").Value = -629
$ws.Range("N132").Value = -20325.845
$ws.Range("H134").Value = 5459.7188
$ws.Range("I134").Value = 4981.357
$ws.Range("J134").Value = 5831.778
$ws.Range("K134").Value = 14944.071
$ws.Range("L134").Value = 17495.334
$ws.Range("M134").Value = -12409.071
$ws.Range("N134").Value = -22565.334
$ws.Range("H136").Value = 5936.091
$ws.Range("I136").Value = 4794.4287
$ws.Range("J136").Value = 6468.8667
$ws.Range("K136").Value = 14383.2861
$ws.Range("L136").Value = 19406.6001
$ws.Range("M136").Value = -11833.2861
$ws.Range("N136").Value = -24506.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 569
$ws.Range("I44").Value = 461.25
$ws.Range("K44").Value = 1383.75
$ws.Range("M44").Value = -985.75
$ws.Range("H82").Value = 8888
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 8888
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 26664
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -27476
$ws.Range("H85").Value = 8888
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 8888
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 26664
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -29472
$ws.Range("H131").Value = 4081.625
$ws.Range("J131").Value = 1674.7241
$ws.Range("L131").Value = 5024.1723
$ws.Range("N131").Value = -15104.1723
$ws.Range("H132").Value = 1317.2
$ws.Range("I132").Value = 986.75
$ws.Range("J132").Value = 1537.5
$ws.Range("K132").Value = 8880.75
$ws.Range("L132").Value = 13837.5
$ws.Range("M132").Value = -6350.75
$ws.Range("N132").Value = -18897.5
$ws.Range("H137").Value = 6968.294
$ws.Range("J137").Value = 8912.091
$ws.Range("L137").Value = 26736.273
$ws.Range("N137").Value = -36936.273
$ws.Range("H139").Value = 4180.7856
$ws.Range("I139").Value = 3803.1
$ws.Range("J139").Value = 5125
$ws.Range("K139").Value = 11409.3
$ws.Range("L139").Value = 15375
$ws.Range("M139").Value = -6269.299999999999
$ws.Range("N139").Value = -25655
$ws.Range("H140").Value = 2048.077
$ws.Range("I140").Value = 685.7143
$ws.Range("J140").Value = 4852.9414
$ws.Range("K140").Value = 2057.1429
$ws.Range("L140").Value = 14558.8242
$ws.Range("M140").Value = 3122.8571
$ws.Range("N140").Value = -24918.8242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 7292086
$ws.Range("J46").Value = 7292086
$ws.Range("L46").Value = 7292086
$ws.Range("N46").Value = -7292398
$ws.Range("H80").Value = 19253.555
$ws.Range("I80").Value = 18000
$ws.Range("K80").Value = 18000
$ws.Range("M80").Value = -17002
$ws.Range("H83").Value = 19253.555
$ws.Range("I83").Value = 18000
$ws.Range("K83").Value = 90000
$ws.Range("M83").Value = -85008
$ws.Range("H107").Value = 527.4666999999999
$ws.Range("I107").Value = 340.10638
$ws.Range("K107").Value = 340.10638
$ws.Range("M107").Value = 1579.89362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5787.706
$ws.Range("I7").Value = 3806.6155
$ws.Range("K7").Value = 3806.6155
$ws.Range("M7").Value = -3694.6155
$ws.Range("H40").Value = 5152.4165
$ws.Range("I40").Value = 5129.3555
$ws.Range("K40").Value = 5129.3555
$ws.Range("M40").Value = -4993.3555
$ws.Range("H46").Value = 11193.333
$ws.Range("I46").Value = 8612.833000000001
$ws.Range("K46").Value = 8612.833000000001
$ws.Range("M46").Value = -8424.833000000001
$ws.Range("H68").Value = 5750.25
$ws.Range("I68").Value = 6499.5
$ws.Range("J68").Value = 5001
$ws.Range("K68").Value = 6499.5
$ws.Range("L68").Value = 5001
$ws.Range("M68").Value = -5750.5
$ws.Range("N68").Value = -6499
$ws.Range("H71").Value = 5750.25
$ws.Range("I71").Value = 6499.5
$ws.Range("J71").Value = 5001
$ws.Range("K71").Value = 32497.5
$ws.Range("L71").Value = 25005
$ws.Range("M71").Value = -28753.5
$ws.Range("N71").Value = -32493
$ws.Range("H82").Value = 100002
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 100002
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H126").Value = 5787.706
$ws.Range("I126").Value = 3806.6155
$ws.Range("K126").Value = 11419.8465
$ws.Range("M126").Value = -8949.8465
$ws.Range("H132").Value = 6141.245
$ws.Range("I132").Value = 6402.5186
$ws.Range("J132").Value = 5820.591
$ws.Range("K132").Value = 19207.5558
$ws.Range("L132").Value = 17461.773
$ws.Range("M132").Value = -16677.5558
$ws.Range("N132").Value = -22521.773

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30025
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H81").Value = 6833.385
$ws.Range("J81").Value = 9248.143
$ws.Range("L81").Value = 18496.286
$ws.Range("N81").Value = -20618.286
$ws.Range("H84").Value = 6833.385
$ws.Range("J84").Value = 9248.143
$ws.Range("L84").Value = 92481.42999999999
$ws.Range("N84").Value = -103089.43
$ws.Range("H110").Value = 8000
$ws.Range("J110").Value = 8000
$ws.Range("L110").Value = 8000
$ws.Range("N110").Value = -16180
$ws.Range("H112").Value = 33281.285
$ws.Range("J112").Value = 33281.285
$ws.Range("L112").Value = 33281.285
$ws.Range("N112").Value = -36235.285
$ws.Range("H126").Value = 12185.667
$ws.Range("I126").Value = 12185.667
$ws.Range("K126").Value = 36557.001
$ws.Range("M126").Value = -34087.001
$ws.Range("H136").Value = 3301.3125
$ws.Range("I136").Value = 2294
$ws.Range("K136").Value = 6882
$ws.Range("M136").Value = -4332
